$wb = $excel.ActiveWorkbook

$ws3 = $wb.Worksheets.Item("Descarga CRM")
$ws4 = $wb.Worksheets.Item("Descarga proyectos")

# --- "Descarga CRM": insert a new column F ("CONTROL PRESUPUESTARIO")
#     between "EQUIPO" (E) and "FECHA ALTA COSTES" (old F) -------------
$ws3.Activate()
$ws3.Columns.Item(6).Insert()

$lo = $ws3.ListObjects.Item(1)
$lo.Resize($ws3.Range("A1:X1048576"))

$ws3.Range("F1").Value = "CONTROL PRESUPUESTARIO"
$ws3.Columns.Item(6).ColumnWidth = 18.25

# Selection / view bookkeeping for "Descarga CRM"
$ws3.Range("F5").Select()

# --- "Descarga proyectos": just a selection / view change -------------
$ws4.Activate()
$ws4.Range("D4").Select()

# --- make "Descarga CRM" the active tab again (matches activeTab=2) ---
$ws3.Activate()
